$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New columns H, I, J for header row (header strings must be entered
#        before the new row-10 label so the shared-string table ends up in
#        the same append order as the target file: H1, I1, then A10, then J1)
$ws.Range("H1").Value = "2022_11_05_09_09_13"
$ws.Range("I1").Value = "2022_11_05_09_14_17"

# --- 2. Insert a brand-new data row above the old "集成结果" row (old row 10)
#        so everything below shifts down by one (old 10->11, 11->12, 12->13).
$ws.Rows("10:10").Insert()

# New row 10 values (new datafile results)
$ws.Range("A10").Value = "2022_11_02_19_58_44-3-0.62829439429"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1

# Header J1 (entered after A10 so shared-string append order matches target)
$ws.Range("J1").Value = "2022_11_05_09_26_44"

# --- 3. Fill in the H/I/J columns for rows 2-9 (existing model rows)
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1

$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1

$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1

$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1

$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 1

# --- 4. Row 11 (the old "集成结果" ensemble-result row, now shifted to 11):
#        fill E11:J11 with newly-computed ensemble scores, and remove the old
#        best-score highlight from B11 (it used to be the yellow-fill "best"
#        cell as B10) since H11 now holds the new best score. Paste the plain
#        (unfilled) format from a neighbouring cell rather than clearing the
#        fill directly - that keeps the shared style table from growing a
#        spurious extra fill entry.
$ws.Range("C11").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("E11").Value = 0.63589020368000004
$ws.Range("F11").Value = 0.63438827984000001
$ws.Range("G11").Value = 0.63079113464000003
$ws.Range("H11").Value = 0.63685043782999995
$ws.Range("I11").Value = 0.62772075951999995
$ws.Range("J11").Value = 0.63410110356000005

# --- 5. Row 12 ("是否修改") and row 13 ("方式") new H/I/J cells
$ws.Range("H12").Value = "是"
$ws.Range("I12").Value = "是"
$ws.Range("J12").Value = "是"

$ws.Range("H13").Value = "权重"
$ws.Range("I13").Value = "权重"
$ws.Range("J13").Value = "权重"

# --- 6. Highlight A5 ("2022_10_25_05_36_40_1-0.62529548458") in red font -
#        added after the row-10 insert so row numbering (A5 unaffected) and
#        new-style ordering line up with the target (font style created
#        before the fill style below).
$ws.Range("A5").Font.Color = 255
$ws.Range("A5").HorizontalAlignment = -4108

# --- 7. New best-score highlight (green fill) on H11, the new max in that row.
$ws.Range("H11").Interior.Color = 5296274
$ws.Range("H11").HorizontalAlignment = -4108

# --- 8. Two brand new trailing rows recording the overall best ensemble run.
#        E15 reuses the same yellow "best value" look as the old B10 cell -
#        set Color before Pattern so the engine matches the existing shared
#        fill instead of minting a near-duplicate one.
$ws.Range("E14").Value = "否"
$ws.Range("E15").Value = 0.64228001062999995
$ws.Range("E15").Interior.Color = 65535
$ws.Range("E15").Interior.Pattern = 1
$ws.Range("E15").HorizontalAlignment = -4108

# --- 9. Column widths for the new H:J columns, matching the bestFit look of
#        the existing E:G columns (E:G already carry the exact original
#        bestFit width, so leave those alone).
$ws.Columns("H:J").ColumnWidth = 19.5

# --- 10. Restore the selection to match the saved worksheet view.
$ws.Range("E24").Select()
